# The commit swaps the bodies of ppt/theme/theme1.xml ("Office Theme") and
# ppt/theme/theme2.xml ("Integral") - file names / relationships are left
# untouched, only the <a:theme> contents trade places. theme2.xml is the
# theme actually bound to the (only) slide master, so it is the part this
# PowerPoint object model can reach; its fontScheme/fmtScheme are already
# byte-identical to theme1's, so the only real difference to reproduce is
# the 12-colour palette (Office Theme's palette moves into theme2.xml).

function ColorRefFromHex([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # COM RGB colours are packed as a Windows COLORREF (0x00BBGGRR)
    return ($b * 65536) + ($g * 256) + $r
}

# Target palette = the "Office Theme" colours that used to live in theme1.xml
$officeThemeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = ColorRefFromHex $officeThemeColors[$i - 1]
}
